$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report row was inserted at row 54 (pushing the former rows
# 54..127 down to 55..128). Insert the row first so everything below it
# shifts down automatically, then populate the newly inserted row.
$ws.Rows("54:54").Insert()

$ws.Range("A54").Value2 = 4
$ws.Range("B54").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value2 = "Los Lagos"
$ws.Range("D54").Value2 = 44467
$ws.Range("E54").Value2 = 10
$ws.Range("F54").Value2 = 100112032
$ws.Range("G54").Value2 = "Zapallo italiano"
$ws.Range("H54").Value2 = "Sin especificar"
$ws.Range("I54").Value2 = "Primera"
$ws.Range("J54").Value2 = 200
$ws.Range("K54").Value2 = 17000
$ws.Range("L54").Value2 = 18000
$ws.Range("M54").Value2 = 17500
$ws.Range("N54").Value2 = "$/caja 50 unidades"
$ws.Range("O54").Value2 = "Región de Arica y Parinacota"
$ws.Range("P54").Value2 = 350
$ws.Range("Q54").Value2 = 50
$ws.Range("R54").Value2 = "Hortaliza"
